$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text to preserve exact original formatting
# (values look numeric but are stored as text; NumberFormat "@" keeps them text,
# then resetting the Style avoids leaving a stray number-format style behind).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.232.77'
$ws.Range('D2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.802.63'
$ws.Range('D3').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.38'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5253'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3822'
$ws.Range('D8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08019'
$ws.Range('D9').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.313'
$ws.Range('D12').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.810.09'
$ws.Range('D15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.314'
$ws.Range('D16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.16'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06605'
$ws.Range('D19').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.36'
$ws.Range('D21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.972'
$ws.Range('D22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.289.41'
$ws.Range('D23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.14'
$ws.Range('D24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.279'
$ws.Range('D25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.70'
$ws.Range('D26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.47'
$ws.Range('D27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.010.50'
$ws.Range('D28').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.29'
$ws.Range('D30').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1083'
$ws.Range('D31').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.057'
$ws.Range('D32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.688'
$ws.Range('D33').Style = 'Normal'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.554'
$ws.Range('D34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07233'
$ws.Range('D35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.40'
$ws.Range('D36').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.118'
$ws.Range('D39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.608'
$ws.Range('D40').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.373'
$ws.Range('D43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.28'
$ws.Range('D44').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.770'
$ws.Range('D46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.24'
$ws.Range('D47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.222'
$ws.Range('D48').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06805'
$ws.Range('D50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.11'
$ws.Range('D51').Style = 'Normal'

# Column E (Volume 1h %) updates - plain text, safe to assign directly
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E8').Value = '  -3.34%  '
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('E16').Value = '  -2.57%  '
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('E18').Value = '  -3.45%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('E26').Value = '  +3.36%  '
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('E30').Value = '  -2.11%  '
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E32').Value = '  -4.45%  '
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('E34').Value = '  -3.98%  '
$ws.Range('E35').Value = '  +3.17%  '
$ws.Range('E36').Value = '  +9.77%  '
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('E38').Value = '  -3.56%  '
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('E41').Value = '  -1.31%  '
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('E47').Value = '  +1.94%  '
$ws.Range('E48').Value = '  +2.93%  '
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('E51').Value = '  -1.63%  '
